# Update "Forecast Comparison" sheet with corrected forecast output:
#  - insert a new "Week_Start_Date" column (B) between "Week" (A) and "ASIN" (now C)
#  - change the Week labels from zero-padded "W01".."W16" to "W1".."W16"
#  - populate the new Week_Start_Date column with the week's start date (as text)

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Forecast Comparison")

# Insert the new column before the current "ASIN" column (B) -> becomes column C.
$ws.Columns("B").Insert()

# Keep the new column formatted as text so dates like "2025-01-05" are stored
# as literal strings, not converted to date serials.
$ws.Columns("B").NumberFormat = "@"

$ws.Range("B1").Value = "Week_Start_Date"

$weeks = @("W1","W2","W3","W4","W5","W6","W7","W8","W9","W10","W11","W12","W13","W14","W15","W16")
$weekStartDates = @(
    "2025-01-05",
    "2025-01-12",
    "2025-01-19",
    "2025-01-26",
    "2025-02-02",
    "2025-02-09",
    "2025-02-16",
    "2025-02-23",
    "2025-03-02",
    "2025-03-09",
    "2025-03-16",
    "2025-03-23",
    "2025-03-30",
    "2025-04-06",
    "2025-04-13",
    "2025-04-20"
)

for ($i = 0; $i -lt $weeks.Length; $i++) {
    $row = $i + 2
    $ws.Range("A$row").Value = $weeks[$i]
    $ws.Range("B$row").Value = $weekStartDates[$i]
}

Write-Output "done"
